$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.192.05'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '1.783.87'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '337.78'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = '0.3914'
$ws.Range('E7').Value = '  +2.09%  '
$ws.Range('D8').Value = '0.3425'
$ws.Range('E8').Value = '  -3.46%  '
$ws.Range('D9').Value = '47.81'
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('D10').Value = '1.189'
$ws.Range('E10').Value = '  -4.01%  '
$ws.Range('D11').Value = '0.07427'
$ws.Range('E11').Value = '  -4.64%  '
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '21.59'
$ws.Range('E13').Value = '  -3.73%  '
$ws.Range('D14').Value = '6.437'
$ws.Range('E14').Value = '  -2.60%  '
$ws.Range('D15').Value = '1.778.19'
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('D16').Value = '7.090'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').Value = '0.00001091'
$ws.Range('E17').Value = '  -2.96%  '
$ws.Range('D18').Value = '0.06655'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').Value = '83.32'
$ws.Range('E19').Value = '  -3.93%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '17.56'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Value = '6.495'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').Value = '27.177.98'
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('D24').Value = '12.34'
$ws.Range('E24').Value = '  -6.48%  '
$ws.Range('D25').Value = '2.374'
$ws.Range('E25').Value = '  -3.70%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '2.503'
$ws.Range('E26').Value = '  -7.19%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '21.16'
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('D28').Value = '1.446'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').Value = '156.00'
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('D31').Value = '134.21'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').Value = '3.976'
$ws.Range('E32').Value = '  -2.34%  '
$ws.Range('D33').Value = '5.990'
$ws.Range('E33').Value = '  -6.09%  '
$ws.Range('D34').Value = '0.08692'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').Value = '13.01'
$ws.Range('E35').Value = '  -6.84%  '
$ws.Range('D36').Value = '1.611'
$ws.Range('E36').Value = '  -4.58%  '
$ws.Range('D37').Value = '5.399'
$ws.Range('E37').Value = '  -4.24%  '
$ws.Range('D38').Value = '0.02390'
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('D39').Value = '0.6787'
$ws.Range('E39').Value = '  -3.78%  '
$ws.Range('D40').Value = '0.06366'
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('D41').Value = '0.2198'
$ws.Range('E41').Value = '  -2.75%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '8.452'
$ws.Range('E42').Value = '  -6.10%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.237'
$ws.Range('E43').Value = '  -4.84%  '
$ws.Range('D44').Value = '14.26'
$ws.Range('E44').Value = '  -3.91%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Value = '0.6382'
$ws.Range('E46').Value = '  -3.69%  '
$ws.Range('D47').Value = '3.857'
$ws.Range('E47').Value = '  -2.74%  '
$ws.Range('D48').Value = '2.133'
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('D49').Value = '131.47'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').Value = '0.07114'
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('D51').Value = '78.66'
$ws.Range('E51').Value = '  -2.81%  '
